# This edit re-shuffles the weekly price records (rows 2-23) of the
# "Fruta, Vega Monumental Concepción - Membrillo" sheet: the contents of
# columns D (Fecha), L..T (Calidad..Kg/unidad) are permuted across rows,
# while the descriptive columns A,B,C,E..K stay untouched (they are
# constant for every row anyway). Row 18 keeps its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together with each logical record.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Maps the destination row (after the edit) to the row that currently
# (before the edit) holds the values that should end up there.
$rowMap = @{
    2  = 10
    3  = 11
    4  = 12
    5  = 7
    6  = 3
    7  = 5
    8  = 6
    9  = 2
    10 = 21
    11 = 23
    12 = 22
    13 = 19
    14 = 16
    15 = 17
    16 = 15
    17 = 20
    18 = 18
    19 = 4
    20 = 8
    21 = 9
    22 = 13
    23 = 14
}

# Snapshot all current values first, since several rows are involved in
# cycles and we must not overwrite a source row before it has been read.
# Note: use Value2 (not Value) - Value's getter is unreliable in this
# COM-interop runtime and returns a description string instead of data.
$snapshot = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 23; $r++) {
        $snapshot["$col$r"] = $ws.Range("$col$r").Value2
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $snapshot["$col$srcRow"]
    }
}
